# Generate Report for Handback
# Updates timestamps / status strings in the handback-status workbook to
# reflect the latest generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) -----------
# 2016-09-05 06:16:30 -> 2016-09-05 06:17:20
$wsOverview.Range("G2").Value = "2016-09-05 06:17:20"
$wsOverview.Range("G4").Value = "2016-09-05 06:17:20"

# --- zh-cn sheet -----------------------------------------------------------
# Priority column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# Correspond Handoff Datetime column (H): 2016-09-05 06:16:26 -> 2016-09-05 06:17:16
$wsZhCn.Range("H2").Value = "2016-09-05 06:17:16"
$wsZhCn.Range("H4").Value = "2016-09-05 06:17:16"

# Correspond Handback DateTime column (K): 2016-09-05 06:16:43 -> 2016-09-05 06:17:33
$wsZhCn.Range("K2").Value = "2016-09-05 06:17:33"
$wsZhCn.Range("K4").Value = "2016-09-05 06:17:33"

# --- de-de sheet -------------------------------------------------------------
# Priority column (E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# Correspond Handoff Datetime column (H): 2016-09-05 06:16:30 -> 2016-09-05 06:17:20
$wsDeDe.Range("H2").Value = "2016-09-05 06:17:20"
$wsDeDe.Range("H4").Value = "2016-09-05 06:17:20"

# Correspond Handback DateTime column (K): 2016-09-05 06:16:50 -> 2016-09-05 06:17:40
$wsDeDe.Range("K2").Value = "2016-09-05 06:17:40"
$wsDeDe.Range("K4").Value = "2016-09-05 06:17:40"
